# Refresh the cryptos list: update Price/Volume(1h) columns for each row,
# and swap in replacement coin rows (34/35 and 48/49) per the upstream feed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D and E hold pre-formatted text (thousand-dot prices, padded
# percentages) in the source workbook; force Text format so COM doesn't
# reinterpret e.g. '147.40' or '7.50' as numbers and drop the trailing zero.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '26.529.70'
$ws.Range("E2").Value = '  -0.41%  '

$ws.Range("D3").Value = '1.627.92'
$ws.Range("E3").Value = '  -0.29%  '

$ws.Range("E4").Value = '  +0.14%  '

$ws.Range("D5").Value = '213.34'
$ws.Range("E5").Value = '  +0.39%  '

$ws.Range("E6").Value = '  +2.27%  '

$ws.Range("E7").Value = '  +0.13%  '

$ws.Range("E8").Value = '  -0.74%  '

$ws.Range("E9").Value = '  +0.06%  '

$ws.Range("D10").Value = '18.79'
$ws.Range("E10").Value = '  -1.09%  '

$ws.Range("E11").Value = '  +0.39%  '

$ws.Range("D12").Value = '1.854.06'
$ws.Range("E12").Value = '  -0.44%  '

$ws.Range("E13").Value = '  +1.68%  '

$ws.Range("D14").Value = '1.604.72'
$ws.Range("E14").Value = '  -2.06%  '

$ws.Range("E15").Value = '  -0.43%  '

$ws.Range("D16").Value = '65.07'
$ws.Range("E16").Value = '  +3.46%  '

$ws.Range("D17").Value = '26.547.81'
$ws.Range("E17").Value = '  -0.39%  '

$ws.Range("E18").Value = '  +0.30%  '

$ws.Range("D19").Value = '214.75'
$ws.Range("E19").Value = '  +2.86%  '

$ws.Range("E20").Value = '  +0.20%  '

$ws.Range("E21").Value = '  -0.14%  '

$ws.Range("D22").Value = '6.25'
$ws.Range("E22").Value = '  +1.32%  '

$ws.Range("E23").Value = '  -0.61%  '

$ws.Range("D24").Value = '2.14'
$ws.Range("E24").Value = '  +11.26%  '

$ws.Range("D25").Value = '147.40'
$ws.Range("E25").Value = '  +0.52%  '

$ws.Range("E26").Value = '  +0.18%  '

$ws.Range("E27").Value = '  +0.29%  '

$ws.Range("E28").Value = '  +2.04%  '

$ws.Range("D29").Value = '15.56'
$ws.Range("E29").Value = '  +1.43%  '

$ws.Range("E30").Value = '  -1.27%  '

$ws.Range("E31").Value = '  -0.53%  '

$ws.Range("E32").Value = '  +3.73%  '

$ws.Range("E33").Value = '  +0.38%  '

$ws.Range("B34").Value = 'Maker'
$ws.Range("C34").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D34").Value = '1.242.19'
$ws.Range("E34").Value = '  +6.36%  '

$ws.Range("B35").Value = 'LidoDAOToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D35").Value = '1.50'
$ws.Range("E35").Value = '  +0.36%  '

$ws.Range("E36").Value = '  +0.33%  '

$ws.Range("E37").Value = '  +4.71%  '

$ws.Range("E38").Value = '  +0.17%  '

$ws.Range("E39").Value = '  +1.45%  '

$ws.Range("D40").Value = '0.795'
$ws.Range("E40").Value = '  -1.46%  '

$ws.Range("E41").Value = '  -2.26%  '

$ws.Range("E42").Value = '  +0.72%  '

$ws.Range("E43").Value = '  -0.83%  '

$ws.Range("D44").Value = '1.764.27'
$ws.Range("E44").Value = '  -0.62%  '

$ws.Range("D45").Value = '93.22'
$ws.Range("E45").Value = '  +1.39%  '

$ws.Range("D46").Value = '1.58'
$ws.Range("E46").Value = '  +2.60%  '

$ws.Range("E47").Value = '  +0.53%  '

$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").Value = '0.0510'
$ws.Range("E48").Value = '  -0.51%  '

$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D49").Value = '0.0958'
$ws.Range("E49").Value = '  +2.34%  '

$ws.Range("D50").Value = '7.50'
$ws.Range("E50").Value = '  -0.17%  '

$ws.Range("E51").Value = '  -0.64%  '
